# [FIX] asset budget report
#
# The "Division" header column is split into two columns:
#   - "Division Code"  (existing column B, renamed)
#   - "Division Name"  (new column, inserted right after, becomes column C)
# Everything that used to live in columns C..S shifts right to D..T.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at C; this pushes the old C:S block to D:T
# and grows the sheet dimension from A1:S9 to A1:T9 automatically.
$ws.Columns("C").Insert()

# Rename the old "Division" header (still in column B) to "Division Code"
# and give the freshly inserted column C its own header, "Division Name".
$ws.Range("B9").Value = "Division Code"
$ws.Range("C9").Value = "Division Name"

# Match the new column widths from the template. The header column (B) is
# now narrower, the new column (C) takes up the space Division used to
# need, and column D (old "Section Code", shifted right) is narrowed to
# match B.
$ws.Columns("B").ColumnWidth = 15.42
$ws.Columns("C").ColumnWidth = 29.45
$ws.Columns("D").ColumnWidth = 15.42

# Selection / scroll position used by the author when saving the sheet.
$ws.Range("A1").Select()
$ws.Range("C5").Select()
